$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.944.94"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.486.12"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.47"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.04"
$ws.Range("E6").Value = "  -4.10%  "
$ws.Range("D7").Value = "3.486.51"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.08"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").Value = "4.085.08"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "3.497.35"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "63.873.82"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.74"
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.59"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "379.01"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").Value = "3.623.14"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.13"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.54"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.17"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("D33").Value = "3.491.01"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.27"
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.25"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.82"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "159.27"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0786"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.87"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.78"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.60"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").Value = "2.411.70"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.886"
$ws.Range("E51").Value = "  -2.61%  "
